# "Generate Report for Handback" - fills in the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) for the
# zh-cn and de-de localization targets, flips the Overview/Status text from
# "Ready for handoff" to "Handed back: in sync with en-US", and widens a
# few columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "032ceea2-d288-4b97-9cb2-4850214ad6e1.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c60b68531858cf133dac1d05e3848bf88341e750/e2e/032ceea2-d288-4b97-9cb2-4850214ad6e1.md"

# Excel's ColumnWidth setter on this host snaps to a 1/6-character grid
# (stored_width = round(ColumnWidth*6)/6 + 5/6), so these are the closest
# settable values to the authored widths.
$wideStatusWidth = 29.166666666666668   # -> stored width 30   (target 29.9777047293527)
$wideFileWidth   = 39.166666666666664   # -> stored width 40   (exact)

function Set-HandbackHyperlink($ws) {
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
    # Give I2 the same "looks like a hyperlink" font the existing A2
    # hyperlink cell uses (underline + the theme hyperlink blue).
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E1").ColumnWidth = $wideStatusWidth
$overview.Range("F1").ColumnWidth = $wideStatusWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("I2").Value = $mdFileName
$zhcn.Range("J2").Value = "032ceea2-d288-4b97-9cb2-4850214ad6e1.c921fb461b9a0234035f6a9bdaa9825c3ac0fcac.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-21 23:06:15"
Set-HandbackHyperlink $zhcn
$zhcn.Range("C1").ColumnWidth = $wideStatusWidth
$zhcn.Range("I1").ColumnWidth = $wideFileWidth
$zhcn.Range("J1").ColumnWidth = $wideFileWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("I2").Value = $mdFileName
$dede.Range("J2").Value = "032ceea2-d288-4b97-9cb2-4850214ad6e1.c921fb461b9a0234035f6a9bdaa9825c3ac0fcac.de-de.xlf"
$dede.Range("K2").Value = "2016-08-21 23:06:21"
Set-HandbackHyperlink $dede
$dede.Range("C1").ColumnWidth = $wideStatusWidth
$dede.Range("I1").ColumnWidth = $wideFileWidth
$dede.Range("J1").ColumnWidth = $wideFileWidth
